# Update the roster table (columns A:C, rows 2-19) with the new
# player / position / team data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Scottie Barnes",     "PG,SG,SF,PF", "Toronto Raptors"),
    @("Josh Giddey",        "PG,SG,SF",    "Chicago Bulls"),
    @("Spencer Dinwiddie",  "PG,SG",       "Dallas Mavericks"),
    @("Nikola Vucevic",     "PF,C",        "Chicago Bulls"),
    @("Evan Mobley",        "PF,C",        "Cleveland Cavaliers"),
    @("Brook Lopez",        "C",           "Milwaukee Bucks"),
    @("Nick Richards",      "C",           "Phoenix Suns"),
    @("Shaedon Sharpe",     "SG,SF",       "Portland Trail Blazers"),
    @("Harrison Barnes",    "SF,PF",       "San Antonio Spurs"),
    @("Mikal Bridges",      "SG,SF,PF",    "New York Knicks"),
    @("De'Aaron Fox",       "PG",          "Sacramento Kings"),
    @("Isaiah Collier",     "PG",          "Utah Jazz"),
    @("Tyler Herro",        "PG,SG",       "Miami Heat"),
    @("Miles Bridges",      "SF,PF",       "Charlotte Hornets"),
    @("DeMar DeRozan",      "SF,PF",       "Sacramento Kings"),
    @("Luka Doncic",        "PG,SG",       "Los Angeles Lakers"),
    @("Bobby Portis",       "PF,C",        "Milwaukee Bucks"),
    @("Ja Morant",          "PG",          "Memphis Grizzlies")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
